$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 updates
$ws.Range("D12").Value = 44539
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 2200
$ws.Range("M12").Value = 2100
$ws.Range("P12").Value = 2100

# Row 13 updates
$ws.Range("D13").Value = 44263
$ws.Range("J13").Value = 140
$ws.Range("K13").Value = 1800
$ws.Range("L13").Value = 2000
$ws.Range("M13").Value = 1914
$ws.Range("N13").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O13").Value = "Provincia de Diguillín"
$ws.Range("P13").Value = 1914
$ws.Range("Q13").Value = 1
